# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310   (the "before" AHB format version)
#   *_new -> *_FV2404   (the "after" AHB format version)
# then wrap the sheet's used range in an Excel Table ("Table1") so the new
# header names are also reflected in the table's column definitions, and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1) Rename the header row (row 1) ------------------------------------
# "Segmentname_old" -> "Segmentname_FV2310", "Segmentname_new" -> "Segmentname_FV2404", etc.
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace "_old$", "_FV2310"
        $newVal = $newVal -replace "_new$", "_FV2404"
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

# --- 2) Turn the data range into a proper Excel Table ---------------------
# (headers are picked up live from row 1, so this automatically uses the
# renamed FV2310 / FV2404 column names)
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header row renamed, Table1 created over" $tableRange.Address() "and header row frozen."
